{"js": "// Replace the date line and each of the 25 division-problem answers in the\n// table with their new values. Every old string in this worksheet is\n// unique, so a simple exact (case-sensitive, whole-match) search/replace\n// per pair is safe and unambiguous.\nconst replacements = [\n    [\"2024-06-22 Saturday\", \"2024-06-23 Sunday\"],\n    [\"832\u00f75=166, 2\", \"349\u00f74=87, 1\"],\n    [\"187\u00f72=93, 1\", \"708\u00f75=141, 3\"],\n    [\"129\u00f72=64, 1\", \"635\u00f76=105, 5\"],\n    [\"743\u00f76=123, 5\", \"913\u00f77=130, 3\"],\n    [\"133\u00f74=33, 1\", \"414\u00f74=103, 2\"],\n    [\"808\u00f78=101, 0\", \"827\u00f78=103, 3\"],\n    [\"993\u00f78=124, 1\", \"417\u00f76=69, 3\"],\n    [\"330\u00f77=47, 1\", \"983\u00f77=140, 3\"],\n    [\"883\u00f73=294, 1\", \"226\u00f77=32, 2\"],\n    [\"701\u00f77=100, 1\", \"152\u00f76=25, 2\"],\n    [\"492\u00f74=123, 0\", \"627\u00f78=78, 3\"],\n    [\"499\u00f77=71, 2\", \"958\u00f77=136, 6\"],\n    [\"674\u00f75=134, 4\", \"150\u00f77=21, 3\"],\n    [\"756\u00f75=151, 1\", \"837\u00f75=167, 2\"],\n    [\"488\u00f78=61, 0\", \"439\u00f74=109, 3\"],\n    [\"588\u00f79=65, 3\", \"576\u00f78=72, 0\"],\n    [\"507\u00f73=169, 0\", \"413\u00f75=82, 3\"],\n    [\"726\u00f79=80, 6\", \"940\u00f79=104, 4\"],\n    [\"878\u00f72=439, 0\", \"103\u00f79=11, 4\"],\n    [\"395\u00f72=197, 1\", \"928\u00f75=185, 3\"],\n    [\"656\u00f76=109, 2\", \"477\u00f76=79, 3\"],\n    [\"832\u00f74=208, 0\", \"855\u00f76=142, 3\"],\n    [\"259\u00f79=28, 7\", \"275\u00f74=68, 3\"],\n    [\"410\u00f74=102, 2\", \"737\u00f75=147, 2\"],\n    [\"594\u00f75=118, 4\", \"759\u00f79=84, 3\"]\n  ];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('text');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error('Could not find text: ' + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each of the 25 division-problem answers in the\n# table with their updated values. Every \"old\" string in this worksheet is\n# unique, so an exact, case-sensitive Find/Replace-All per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '2024-06-22 Saturday'; New = '2024-06-23 Sunday' },\n    @{ Old = '832\u00f75=166, 2'; New = '349\u00f74=87, 1' },\n    @{ Old = '187\u00f72=93, 1'; New = '708\u00f75=141, 3' },\n    @{ Old = '129\u00f72=64, 1'; New = '635\u00f76=105, 5' },\n    @{ Old = '743\u00f76=123, 5'; New = '913\u00f77=130, 3' },\n    @{ Old = '133\u00f74=33, 1'; New = '414\u00f74=103, 2' },\n    @{ Old = '808\u00f78=101, 0'; New = '827\u00f78=103, 3' },\n    @{ Old = '993\u00f78=124, 1'; New = '417\u00f76=69, 3' },\n    @{ Old = '330\u00f77=47, 1'; New = '983\u00f77=140, 3' },\n    @{ Old = '883\u00f73=294, 1'; New = '226\u00f77=32, 2' },\n    @{ Old = '701\u00f77=100, 1'; New = '152\u00f76=25, 2' },\n    @{ Old = '492\u00f74=123, 0'; New = '627\u00f78=78, 3' },\n    @{ Old = '499\u00f77=71, 2'; New = '958\u00f77=136, 6' },\n    @{ Old = '674\u00f75=134, 4'; New = '150\u00f77=21, 3' },\n    @{ Old = '756\u00f75=151, 1'; New = '837\u00f75=167, 2' },\n    @{ Old = '488\u00f78=61, 0'; New = '439\u00f74=109, 3' },\n    @{ Old = '588\u00f79=65, 3'; New = '576\u00f78=72, 0' },\n    @{ Old = '507\u00f73=169, 0'; New = '413\u00f75=82, 3' },\n    @{ Old = '726\u00f79=80, 6'; New = '940\u00f79=104, 4' },\n    @{ Old = '878\u00f72=439, 0'; New = '103\u00f79=11, 4' },\n    @{ Old = '395\u00f72=197, 1'; New = '928\u00f75=185, 3' },\n    @{ Old = '656\u00f76=109, 2'; New = '477\u00f76=79, 3' },\n    @{ Old = '832\u00f74=208, 0'; New = '855\u00f76=142, 3' },\n    @{ Old = '259\u00f79=28, 7'; New = '275\u00f74=68, 3' },\n    @{ Old = '410\u00f74=102, 2'; New = '737\u00f75=147, 2' },\n    @{ Old = '594\u00f75=118, 4'; New = '759\u00f79=84, 3' }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # wdFindContinue=1, wdReplaceAll=2\n    $found = $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        throw \"Could not find text: $($pair.Old)\"\n    }\n}\n"}
